$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 92.59999999999999
$ws.Range("I11").Value = 92.59999999999999
$ws.Range("K11").Value = 92.59999999999999
$ws.Range("M11").Value = 47.40000000000001
$ws.Range("H28").Value = 952.55554
$ws.Range("I28").Value = 920.35297
$ws.Range("K28").Value = 920.35297
$ws.Range("M28").Value = -435.35297
$ws.Range("H43").Value = 3162.6667
$ws.Range("I43").Value = 3999
$ws.Range("J43").Value = 2744.5
$ws.Range("K43").Value = 3999
$ws.Range("L43").Value = 2744.5
$ws.Range("M43").Value = -3930
$ws.Range("N43").Value = -2882.5
$ws.Range("H53").Value = 5897.4
$ws.Range("J53").Value = 814.3333
$ws.Range("L53").Value = 814.3333
$ws.Range("N53").Value = -2088.3333
$ws.Range("H103").Value = 1280.5
$ws.Range("I103").Value = 280.625
$ws.Range("J103").Value = 2080.4
$ws.Range("K103").Value = 841.875
$ws.Range("L103").Value = 6241.200000000001
$ws.Range("M103").Value = -255.875
$ws.Range("N103").Value = -7413.200000000001
$ws.Range("H112").Value = 1718.4736
$ws.Range("I112").Value = 1395.6
$ws.Range("J112").Value = 1833.7858
$ws.Range("K112").Value = 4186.799999999999
$ws.Range("L112").Value = 5501.357400000001
$ws.Range("M112").Value = -3078.799999999999
$ws.Range("N112").Value = -7717.357400000001
$ws.Range("H125").Value = 2974.25
$ws.Range("I125").Value = 3185.889
$ws.Range("J125").Value = 2801.0908
$ws.Range("K125").Value = 28673.001
$ws.Range("L125").Value = 25209.8172
$ws.Range("M125").Value = -26213.001
$ws.Range("N125").Value = -30129.8172
$ws.Range("H132").Value = 3624.3333
$ws.Range("I132").Value = 3649.6667
$ws.Range("K132").Value = 10949.0001
$ws.Range("M132").Value = -8419.000100000001
$ws.Range("H133").Value = 147721.38
$ws.Range("J133").Value = 147721.38
$ws.Range("L133").Value = 147721.38
$ws.Range("N133").Value = -157841.38
$ws.Range("H138").Value = 4621.3335
$ws.Range("I138").Value = 3133
$ws.Range("K138").Value = 9399
$ws.Range("M138").Value = -4259
$ws.Range("H139").Value = 96744.664
$ws.Range("J139").Value = 96744.664
$ws.Range("L139").Value = 96744.664
$ws.Range("N139").Value = -107024.664
$ws.Range("H140").Value = 69279.47
$ws.Range("J140").Value = 69279.47
$ws.Range("L140").Value = 69279.47
$ws.Range("N140").Value = -79639.47

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1179.625
$ws.Range("I63").Value = 1143.75
$ws.Range("J63").Value = 1287.25
$ws.Range("K63").Value = 1143.75
$ws.Range("L63").Value = 1287.25
$ws.Range("M63").Value = -457.75
$ws.Range("N63").Value = -2659.25
$ws.Range("H64").Value = 59000
$ws.Range("J64").Value = 59000
$ws.Range("L64").Value = 59000
$ws.Range("N64").Value = -59496
$ws.Range("H66").Value = 1179.625
$ws.Range("I66").Value = 1143.75
$ws.Range("J66").Value = 1287.25
$ws.Range("K66").Value = 5718.75
$ws.Range("L66").Value = 6436.25
$ws.Range("M66").Value = -2286.75
$ws.Range("N66").Value = -13300.25
$ws.Range("H67").Value = 59000
$ws.Range("J67").Value = 59000
$ws.Range("L67").Value = 59000
$ws.Range("N67").Value = -60716
$ws.Range("H74").Value = 4304.231
$ws.Range("I74").Value = 2368.4
$ws.Range("K74").Value = 2368.4
$ws.Range("M74").Value = -1494.4
$ws.Range("H77").Value = 4304.231
$ws.Range("I77").Value = 2368.4
$ws.Range("K77").Value = 11842
$ws.Range("M77").Value = -7474
$ws.Range("H108").Value = 69380
$ws.Range("J108").Value = 69380
$ws.Range("L108").Value = 69380
$ws.Range("N108").Value = -77060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 808.5333000000001
$ws.Range("I22").Value = 802.2308
$ws.Range("J22").Value = 849.5
$ws.Range("K22").Value = 802.2308
$ws.Range("L22").Value = 849.5
$ws.Range("M22").Value = -629.2308
$ws.Range("N22").Value = -1195.5
$ws.Range("H86").Value = 111113620
$ws.Range("I86").Value = 333334940
$ws.Range("J86").Value = 2963.6667
$ws.Range("K86").Value = 333334940
$ws.Range("L86").Value = 2963.6667
$ws.Range("M86").Value = -333333817
$ws.Range("N86").Value = -5209.6667
$ws.Range("H89").Value = 111113620
$ws.Range("I89").Value = 333334940
$ws.Range("J89").Value = 2963.6667
$ws.Range("K89").Value = 1666674700
$ws.Range("L89").Value = 14818.3335
$ws.Range("M89").Value = -1666669084
$ws.Range("N89").Value = -26050.3335
$ws.Range("H105").Value = 125001570
$ws.Range("I105").Value = 125001570
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 125001570
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -124999823

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 5175.75
$ws.Range("I5").Value = 6737.6665
$ws.Range("J5").Value = 490
$ws.Range("K5").Value = 6737.6665
$ws.Range("L5").Value = 490
$ws.Range("M5").Value = -6625.6665
$ws.Range("N5").Value = -714
$ws.Range("H25").Value = 22476.5
$ws.Range("I25").Value = 23971.8
$ws.Range("J25").Value = 15000
$ws.Range("K25").Value = 23971.8
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = -23797.8
$ws.Range("N25").Value = -15348
$ws.Range("H99").Value = 2919
$ws.Range("I99").Value = 2878.6
$ws.Range("J99").Value = 2959.4
$ws.Range("K99").Value = 2878.6
$ws.Range("L99").Value = 2959.4
$ws.Range("M99").Value = -1380.6
$ws.Range("N99").Value = -5955.4
$ws.Range("H107").Value = 2838.6191
$ws.Range("I107").Value = 2397.8
$ws.Range("J107").Value = 3940.6667
$ws.Range("K107").Value = 2397.8
$ws.Range("L107").Value = 3940.6667
$ws.Range("M107").Value = -477.8000000000002
$ws.Range("N107").Value = -7780.6667
$ws.Range("H126").Value = 2919
$ws.Range("I126").Value = 2878.6
$ws.Range("J126").Value = 2959.4
$ws.Range("K126").Value = 8635.799999999999
$ws.Range("L126").Value = 8878.200000000001
$ws.Range("M126").Value = -6165.799999999999
$ws.Range("N126").Value = -13818.2
$ws.Range("H131").Value = 67500
$ws.Range("J131").Value = 67500
$ws.Range("L131").Value = 67500
$ws.Range("N131").Value = -77580
$ws.Range("H138").Value = 73076.78
$ws.Range("J138").Value = 73076.78
$ws.Range("L138").Value = 73076.78
$ws.Range("N138").Value = -83356.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 554919.5600000001
$ws.Range("J4").Value = 3674
$ws.Range("L4").Value = 11022
$ws.Range("N4").Value = -11246
$ws.Range("H13").Value = 313.125
$ws.Range("I13").Value = 356.7143
$ws.Range("J13").Value = 8
$ws.Range("K13").Value = 1070.1429
$ws.Range("L13").Value = 24
$ws.Range("M13").Value = -902.1428999999998
$ws.Range("N13").Value = -360
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1350.1666
$ws.Range("K46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("M46").Value = 4050.4998
$ws.Range("N46").Value = -4232.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5651.5713
$ws.Range("I70").Value = 4903.25
$ws.Range("J70").Value = 6649.3335
$ws.Range("K70").Value = 4903.25
$ws.Range("L70").Value = 6649.3335
$ws.Range("M70").Value = -4633.25
$ws.Range("N70").Value = -7189.3335
$ws.Range("H73").Value = 5651.5713
$ws.Range("I73").Value = 4903.25
$ws.Range("J73").Value = 6649.3335
$ws.Range("K73").Value = 4903.25
$ws.Range("L73").Value = 6649.3335
$ws.Range("M73").Value = -3967.25
$ws.Range("N73").Value = -8521.333500000001
$ws.Range("H126").Value = 25004094
$ws.Range("I126").Value = 55558660
$ws.Range("J126").Value = 4902
$ws.Range("K126").Value = 166675980
$ws.Range("L126").Value = 14706
$ws.Range("M126").Value = -166673510
$ws.Range("N126").Value = -19646

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2959.5588
$ws.Range("I40").Value = 2599.2593
$ws.Range("J40").Value = 4349.2856
$ws.Range("K40").Value = 2599.2593
$ws.Range("L40").Value = 4349.2856
$ws.Range("M40").Value = -2463.2593
$ws.Range("N40").Value = -4621.2856
$ws.Range("H122").Value = 4323
$ws.Range("I122").Value = 3256.5173
$ws.Range("K122").Value = 9769.5519
$ws.Range("M122").Value = -7319.5519
$ws.Range("H132").Value = 16480.732
$ws.Range("I132").Value = 14443.643
$ws.Range("K132").Value = 43330.929
$ws.Range("M132").Value = -40800.929

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8768.200000000001
$ws.Range("I62").Value = 4996
$ws.Range("J62").Value = 10799.385
$ws.Range("K62").Value = 4996
$ws.Range("L62").Value = 10799.385
$ws.Range("M62").Value = -4372
$ws.Range("N62").Value = -12047.385
$ws.Range("H65").Value = 8768.200000000001
$ws.Range("I65").Value = 4996
$ws.Range("J65").Value = 10799.385
$ws.Range("K65").Value = 24980
$ws.Range("L65").Value = 53996.925
$ws.Range("M65").Value = -21860
$ws.Range("N65").Value = -60236.925
$ws.Range("H126").Value = 2601.6667
$ws.Range("I126").Value = 2432.2
$ws.Range("K126").Value = 7296.599999999999
$ws.Range("M126").Value = -4826.599999999999
$ws.Range("H135").Value = 72000
$ws.Range("J135").Value = 72000
$ws.Range("L135").Value = 72000
$ws.Range("N135").Value = -82140
